$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 53, shifting existing rows 53:69 down to 54:70
$ws.Rows(53).Insert()

# Populate the new row 53 with the new weekly record
$ws.Cells.Item(53, 1).Value = 6
$ws.Cells.Item(53, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(53, 3).Value = "Metropolitana"
$ws.Cells.Item(53, 4).Value = 44736
$ws.Cells.Item(53, 5).Value = 13
$ws.Cells.Item(53, 6).Value = "Fruta"
$ws.Cells.Item(53, 7).Value = 100108
$ws.Cells.Item(53, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(53, 9).Value = 100108007
$ws.Cells.Item(53, 10).Value = "Coco"
$ws.Cells.Item(53, 11).Value = "Sin especificar"
$ws.Cells.Item(53, 12).Value = "Primera"
$ws.Cells.Item(53, 13).Value = 100
$ws.Cells.Item(53, 14).Value = 22000
$ws.Cells.Item(53, 15).Value = 22000
$ws.Cells.Item(53, 16).Value = 22000
$ws.Cells.Item(53, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(53, 18).Value = "Perú"
$ws.Cells.Item(53, 19).Value = 1100
$ws.Cells.Item(53, 20).Value = 20
